$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.828.97'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +13.28%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.731.90'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +7.40%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9968'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.89%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.39'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.34%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9956'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.58%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3789'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +4.04%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3623'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +7.11%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '50.40'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +22.02%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.222'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +8.07%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07621'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +8.78%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9991'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.34%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.67'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +9.92%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.456'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +10.43%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.082'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +7.79%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.733.21'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +8.48%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001149'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +7.03%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9941'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.26%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06821'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.27%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '86.38'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +11.95%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.37'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +8.52%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.408'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +7.71%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.64'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +9.05%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '25.642.71'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +13.21%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.445'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.51%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.915'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +14.74%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.49'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +6.07%  '

$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '154.58'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.15%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '135.47'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +9.23%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.924.89'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +8.75%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.181'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +24.55%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.929'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +15.87%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.107'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.89%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.814'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +10.81%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '13.78'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +15.10%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08643'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +6.02%  '

$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06724'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +11.61%  '

$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.576'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +8.89%  '

$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '9.229'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +8.29%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02460'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +11.98%  '

$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2203'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +10.29%  '

$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.288'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.77%  '

$ws.Range('E43').Value = '  +10.17%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9946'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.32%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.78'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +8.12%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6256'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +10.35%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.882'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.55%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.138'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +9.02%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '131.58'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +5.50%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07453'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +8.71%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.07'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +8.10%  '
